# Auto-generated script to update Recommandations and Top_YTD sheets
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Recommandations")

$arr1 = New-Object 'object[,]' 44,7
$arr1[0,0] = "BRVM - SERVICES PUBLICS"
$arr1[0,1] = 0
$arr1[0,2] = 8
$arr1[0,3] = 3372.77
$arr1[0,4] = 107.66
$arr1[0,5] = "🟡 Observer"
$arr1[0,6] = "➖ Neutre"
$arr1[1,0] = "NEI-CEDA CI"
$arr1[1,1] = 0
$arr1[1,2] = 4
$arr1[1,3] = 2735
$arr1[1,4] = 670
$arr1[1,5] = "🟡 Observer"
$arr1[1,6] = "➖ Neutre"
$arr1[2,0] = "AIR LIQUIDE CI"
$arr1[2,1] = 0
$arr1[2,2] = 4
$arr1[2,3] = 2685
$arr1[2,4] = 670
$arr1[2,5] = "🟡 Observer"
$arr1[2,6] = "➖ Neutre"
$arr1[3,0] = "BRVM - AUTRES SECTEURS"
$arr1[3,1] = 0
$arr1[3,2] = 4
$arr1[3,3] = 2472.08
$arr1[3,4] = 625.47
$arr1[3,5] = "🟡 Observer"
$arr1[3,6] = "➖ Neutre"
$arr1[4,0] = "BRVM - DISTRIBUTION"
$arr1[4,1] = 0
$arr1[4,2] = 4
$arr1[4,3] = 2057.32
$arr1[4,4] = 533.53
$arr1[4,5] = "🟡 Observer"
$arr1[4,6] = "➖ Neutre"
$arr1[5,0] = "BRVM - AGRICULTURE"
$arr1[5,1] = 0
$arr1[5,2] = 4
$arr1[5,3] = 1489.93
$arr1[5,4] = 374.65
$arr1[5,5] = "🟡 Observer"
$arr1[5,6] = "➖ Neutre"
$arr1[6,0] = "BRVM - TRANSPORT"
$arr1[6,1] = 0
$arr1[6,2] = 4
$arr1[6,3] = 1445.38
$arr1[6,4] = 353.7
$arr1[6,5] = "🟡 Observer"
$arr1[6,6] = "➖ Neutre"
$arr1[7,0] = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$arr1[7,1] = 0
$arr1[7,2] = 4
$arr1[7,3] = 720.42
$arr1[7,4] = 193.87
$arr1[7,5] = "🟡 Observer"
$arr1[7,6] = "➖ Neutre"
$arr1[8,0] = "BRVM-PRESTIGE"
$arr1[8,1] = 0
$arr1[8,2] = 4
$arr1[8,3] = 557.11
$arr1[8,4] = 138.63
$arr1[8,5] = "🟡 Observer"
$arr1[8,6] = "➖ Neutre"
$arr1[9,0] = "BRVM - FINANCES"
$arr1[9,1] = 0
$arr1[9,2] = 4
$arr1[9,3] = 556.0599999999999
$arr1[9,4] = 139.34
$arr1[9,5] = "🟡 Observer"
$arr1[9,6] = "➖ Neutre"
$arr1[10,0] = "BRVM - SERVICES FINANCIERS"
$arr1[10,1] = 0
$arr1[10,2] = 4
$arr1[10,3] = 546.48
$arr1[10,4] = 136.94
$arr1[10,5] = "🟡 Observer"
$arr1[10,6] = "➖ Neutre"
$arr1[11,0] = "BRVM - INDUSTRIELS"
$arr1[11,1] = 0
$arr1[11,2] = 4
$arr1[11,3] = 518.14
$arr1[11,4] = 130.44
$arr1[11,5] = "🟡 Observer"
$arr1[11,6] = "➖ Neutre"
$arr1[12,0] = "BRVM - ENERGIE"
$arr1[12,1] = 0
$arr1[12,2] = 4
$arr1[12,3] = 444.91
$arr1[12,4] = 108.71
$arr1[12,5] = "🟡 Observer"
$arr1[12,6] = "➖ Neutre"
$arr1[13,0] = "BRVM - TELECOMMUNICATIONS"
$arr1[13,1] = 0
$arr1[13,2] = 4
$arr1[13,3] = 381.73
$arr1[13,4] = 94.77
$arr1[13,5] = "🟡 Observer"
$arr1[13,6] = "➖ Neutre"
$arr1[14,0] = "BRVM - INDUSTRIE              (**)"
$arr1[14,1] = 0
$arr1[14,2] = 1
$arr1[14,3] = 214.53
$arr1[14,4] = 214.53
$arr1[14,5] = "🟡 Observer"
$arr1[14,6] = "➖ Neutre"
$arr1[15,0] = "BRVM-PRINCIPAL                (**)"
$arr1[15,1] = 0
$arr1[15,2] = 1
$arr1[15,3] = 205.03
$arr1[15,4] = 205.03
$arr1[15,5] = "🟡 Observer"
$arr1[15,6] = "➖ Neutre"
$arr1[16,0] = "BRVM - CONSOMMATION DE BASE               (**)"
$arr1[16,1] = 0
$arr1[16,2] = 1
$arr1[16,3] = 191.47
$arr1[16,4] = 191.47
$arr1[16,5] = "🟡 Observer"
$arr1[16,6] = "➖ Neutre"
$arr1[17,0] = "SAFCA CI (SAFC)"
$arr1[17,1] = 4
$arr1[17,2] = 0
$arr1[17,3] = 29.54
$arr1[17,4] = 7.36
$arr1[17,5] = "🟢 Achat"
$arr1[17,6] = "✅ Renforcer"
$arr1[18,0] = "SICABLE CI (CABC)"
$arr1[18,1] = 3
$arr1[18,2] = 0
$arr1[18,3] = 22.02
$arr1[18,4] = 7.36
$arr1[18,5] = "🟢 Achat"
$arr1[18,6] = "✅ Renforcer"
$arr1[19,0] = "CFAO MOTORS CI (CFAC)"
$arr1[19,1] = 3
$arr1[19,2] = 0
$arr1[19,3] = 20.28
$arr1[19,4] = 7.35
$arr1[19,5] = "🟢 Achat"
$arr1[19,6] = "✅ Renforcer"
$arr1[20,0] = "UNIWAX CI (UNXC)"
$arr1[20,1] = 2
$arr1[20,2] = 0
$arr1[20,3] = 12.15
$arr1[20,4] = 7.42
$arr1[20,5] = "🟡 Observer"
$arr1[20,6] = "➖ Neutre"
$arr1[21,0] = "SERVAIR ABIDJAN CI (ABJC)"
$arr1[21,1] = 1
$arr1[21,2] = 0
$arr1[21,3] = 7.38
$arr1[21,4] = 7.38
$arr1[21,5] = "🟡 Observer"
$arr1[21,6] = "➖ Neutre"
$arr1[22,0] = "BICI CI (BICC)"
$arr1[22,1] = 1
$arr1[22,2] = 0
$arr1[22,3] = 5.26
$arr1[22,4] = 5.26
$arr1[22,5] = "🟡 Observer"
$arr1[22,6] = "➖ Neutre"
$arr1[23,0] = "FILTISAC CI (FTSC)"
$arr1[23,1] = 2
$arr1[23,2] = 2
$arr1[23,3] = 5.11
$arr1[23,4] = -3.07
$arr1[23,5] = "🟡 Observer"
$arr1[23,6] = "👀 À surveiller"
$arr1[24,0] = "BANK OF AFRICA SENEGAL (BOAS)"
$arr1[24,1] = 1
$arr1[24,2] = 0
$arr1[24,3] = 4.6
$arr1[24,4] = 4.6
$arr1[24,5] = "🟡 Observer"
$arr1[24,6] = "➖ Neutre"
$arr1[25,0] = "LOTERIE NATIONALE DU BENIN (LNBB)"
$arr1[25,1] = 1
$arr1[25,2] = 0
$arr1[25,3] = 4.08
$arr1[25,4] = 4.08
$arr1[25,5] = "🟡 Observer"
$arr1[25,6] = "➖ Neutre"
$arr1[26,0] = "NSIA BANQUE COTE D'IVOIRE (NSBC)"
$arr1[26,1] = 1
$arr1[26,2] = 0
$arr1[26,3] = 2.23
$arr1[26,4] = 2.23
$arr1[26,5] = "🟡 Observer"
$arr1[26,6] = "➖ Neutre"
$arr1[27,0] = "ECOBANK TRANS. INCORP. TG (ETIT)"
$arr1[27,1] = 1
$arr1[27,2] = 1
$arr1[27,3] = 0.26
$arr1[27,4] = -5
$arr1[27,5] = "🟡 Observer"
$arr1[27,6] = "👀 À surveiller"
$arr1[28,0] = "TOTAL"
$arr1[28,1] = 0
$arr1[28,2] = 3
$arr1[28,3] = 0
$arr1[28,4] = 0
$arr1[28,5] = "🟡 Observer"
$arr1[28,6] = "➖ Neutre"
$arr1[29,0] = "SAPH CI (SPHC)"
$arr1[29,1] = 0
$arr1[29,2] = 1
$arr1[29,3] = -1.31
$arr1[29,4] = -1.31
$arr1[29,5] = "🟡 Observer"
$arr1[29,6] = "➖ Neutre"
$arr1[30,0] = "PALM CI (PALC)"
$arr1[30,1] = 0
$arr1[30,2] = 1
$arr1[30,3] = -1.45
$arr1[30,4] = -1.45
$arr1[30,5] = "🟡 Observer"
$arr1[30,6] = "➖ Neutre"
$arr1[31,0] = "SOGB CI (SOGC)"
$arr1[31,1] = 0
$arr1[31,2] = 1
$arr1[31,3] = -1.58
$arr1[31,4] = -1.58
$arr1[31,5] = "🟡 Observer"
$arr1[31,6] = "➖ Neutre"
$arr1[32,0] = "ORAGROUP TOGO (ORGT)"
$arr1[32,1] = 0
$arr1[32,2] = 1
$arr1[32,3] = -1.6
$arr1[32,4] = -1.6
$arr1[32,5] = "🟡 Observer"
$arr1[32,6] = "➖ Neutre"
$arr1[33,0] = "AIR LIQUIDE CI (SIVC)"
$arr1[33,1] = 0
$arr1[33,2] = 1
$arr1[33,3] = -2.24
$arr1[33,4] = -2.24
$arr1[33,5] = "🟡 Observer"
$arr1[33,6] = "➖ Neutre"
$arr1[34,0] = "NEI-CEDA CI (NEIC)"
$arr1[34,1] = 0
$arr1[34,2] = 1
$arr1[34,3] = -2.9
$arr1[34,4] = -2.9
$arr1[34,5] = "🟡 Observer"
$arr1[34,6] = "➖ Neutre"
$arr1[35,0] = "BANK OF AFRICA NG (BOAN)"
$arr1[35,1] = 0
$arr1[35,2] = 1
$arr1[35,3] = -3.7
$arr1[35,4] = -3.7
$arr1[35,5] = "🟡 Observer"
$arr1[35,6] = "➖ Neutre"
$arr1[36,0] = "ONATEL BF (ONTBF)"
$arr1[36,1] = 0
$arr1[36,2] = 1
$arr1[36,3] = -3.85
$arr1[36,4] = -3.85
$arr1[36,5] = "🟡 Observer"
$arr1[36,6] = "➖ Neutre"
$arr1[37,0] = "SOCIETE GENERALE COTE D'IVOIRE (SGBC)"
$arr1[37,1] = 0
$arr1[37,2] = 1
$arr1[37,3] = -3.91
$arr1[37,4] = -3.91
$arr1[37,5] = "🟡 Observer"
$arr1[37,6] = "➖ Neutre"
$arr1[38,0] = "TRACTAFRIC MOTORS CI (PRSC)"
$arr1[38,1] = 0
$arr1[38,2] = 2
$arr1[38,3] = -4.21
$arr1[38,4] = -1.43
$arr1[38,5] = "🟡 Observer"
$arr1[38,6] = "➖ Neutre"
$arr1[39,0] = "CIE CI (CIEC)"
$arr1[39,1] = 0
$arr1[39,2] = 1
$arr1[39,3] = -6.3
$arr1[39,4] = -6.3
$arr1[39,5] = "🟡 Observer"
$arr1[39,6] = "➖ Neutre"
$arr1[40,0] = "SICOR CI (SICC)"
$arr1[40,1] = 0
$arr1[40,2] = 1
$arr1[40,3] = -6.91
$arr1[40,4] = -6.91
$arr1[40,5] = "🟡 Observer"
$arr1[40,6] = "➖ Neutre"
$arr1[41,0] = "TOTALENERGIES MARKETING CI (TTLC)"
$arr1[41,1] = 0
$arr1[41,2] = 1
$arr1[41,3] = -7
$arr1[41,4] = -7
$arr1[41,5] = "🟡 Observer"
$arr1[41,6] = "➖ Neutre"
$arr1[42,0] = "SETAO CI (STAC)"
$arr1[42,1] = 0
$arr1[42,2] = 1
$arr1[42,3] = -7.2
$arr1[42,4] = -7.2
$arr1[42,5] = "🟡 Observer"
$arr1[42,6] = "➖ Neutre"
$arr1[43,0] = "SOLIBRA CI (SLBC)"
$arr1[43,1] = 0
$arr1[43,2] = 2
$arr1[43,3] = -7.5
$arr1[43,4] = -3.78
$arr1[43,5] = "🟡 Observer"
$arr1[43,6] = "➖ Neutre"
$ws1.Range("A2:G45").Value = $arr1

$ws2 = $wb.Worksheets.Item("Top_YTD")
$arr2 = New-Object 'object[,]' 10,1
$arr2[0,0] = 9424660.609999999
$arr2[1,0] = 377137.25
$arr2[2,0] = 353698.16
$arr2[3,0] = 265664.43
$arr2[4,0] = 142166.47
$arr2[5,0] = 49732.05
$arr2[6,0] = 45191.75
$arr2[7,0] = 6042.09
$arr2[8,0] = 3177.96
$arr2[9,0] = 3163.59
$ws2.Range("B2:B11").Value = $arr2

